$d = $word.ActiveDocument

# Pull the whole package as flattened WordOpenXML so we can make precise,
# surgical text edits to word/styles.xml (element removal/addition isn't
# reachable purely through the Style/Font/ParagraphFormat property setters).
$xml = $d.WordOpenXML

# --- Title / TitleChar run properties: drop the condensed-spacing / kerning
#     pair that used to sit between rFonts and sz in both style defs.
$xml = $xml.Replace('<w:spacing w:val="-10"/><w:kern w:val="28"/><w:sz w:val="56"/><w:szCs w:val="56"/></w:rPr>', '<w:sz w:val="56"/><w:szCs w:val="56"/></w:rPr>')

# --- Author style: base it on Title, drop the explicit center alignment
#     (now inherited from Title), add an explicit run size override.
$authorOld = '<w:style w:type="paragraph" w:customStyle="1" w:styleId="Author"><w:name w:val="Author"/><w:next w:val="BodyText"/><w:qFormat/><w:pPr><w:keepNext/><w:keepLines/><w:jc w:val="center"/></w:pPr></w:style>'
$authorNew = '<w:style w:type="paragraph" w:customStyle="1" w:styleId="Author"><w:name w:val="Author"/><w:basedOn w:val="Title"/><w:next w:val="BodyText"/><w:qFormat/><w:pPr><w:keepNext/><w:keepLines/></w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:style>'
$xml = $xml.Replace($authorOld, $authorNew)

# --- Date style: same treatment as Author.
$dateOld = '<w:style w:type="paragraph" w:styleId="Date"><w:name w:val="Date"/><w:next w:val="BodyText"/><w:qFormat/><w:pPr><w:keepNext/><w:keepLines/><w:jc w:val="center"/></w:pPr></w:style>'
$dateNew = '<w:style w:type="paragraph" w:styleId="Date"><w:name w:val="Date"/><w:basedOn w:val="Title"/><w:next w:val="BodyText"/><w:qFormat/><w:pPr><w:keepNext/><w:keepLines/></w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:style>'
$xml = $xml.Replace($dateOld, $dateNew)

$d.WordOpenXML = $xml
